# Applies the diary update described in the commit "Updated diary as of 6 Feb 2020 (#261)"
# Fills in rows 28-31 of the diary worksheet with new entries and updates the
# sheet's scroll/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: 30 Jan 2020 (Th) entry continued
$ws.Range("A28").Value = "30 Jan 2020 (Th)"
$ws.Range("B28").Value = "1545-1700"
$ws.Range("C28").Value = "Harry, Deon, Thuc"
$ws.Range("D28").Value = "Finalize writeup, finalize UML diagram"
$ws.Range("E28").Value = "Mission accomplished"

# Row 29: 5 Jan 2020 (W) entry
$ws.Range("A29").Value = "5 Jan 2020 (W)"
$ws.Range("B29").Value = "1710-2010"
$ws.Range("C29").Value = "Harry, Deon, Thuc"
$ws.Range("D29").Value = "Work on fourth lecture’s homework while simutaneously paying attention to the 261 lecture"
$ws.Range("E29").Value = "We finished documenting the first feature"
$ws.Range("F29").Value = "Decisions, decisions! It’s hard to make a choice of features when there are so many! Also in Runeline a large chunk of the code is plugins. So the hunt for essential features became a lot harder."
$ws.Range("G29").Value = "Because the classes are in sequence, I have to work on the assignments in sequence. So the work for this class will always be at the tail end of the weekly schedule."

# Row 30: continuation of 5 Jan 2020 (W), some cells reference row 29 via formula
$ws.Range("A30").Value = "5 Jan 2020 (W)"
$ws.Range("B30").Value = "2022-0000"
$ws.Range("C30").Formula = "=C29"
$ws.Range("D30").Formula = "=D29"
$ws.Range("E30").Value = "We finished documenting the second feature"

# Row 31: 6 Jan 2020 (Th) entry
$ws.Range("A31").Value = "6 Jan 2020 (Th)"
$ws.Range("B31").Value = "0000-0030"
$ws.Range("C31").Value = "Harry, Deon, Thuc"
$ws.Range("D31").Formula = "=D30"
$ws.Range("E31").Formula = "=E30"

# Row heights as observed in the final workbook
$ws.Rows.Item(28).RowHeight = 25.35
$ws.Rows.Item(29).RowHeight = 73.1
$ws.Rows.Item(30).RowHeight = 37.3
$ws.Rows.Item(31).RowHeight = 37.3

# Update the view's scroll position / selection to match the edited region
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
